$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row 40
$ws.Range("H40").Value = 1558.6471
$ws.Range("J40").Value = 1508.4166
$ws.Range("L40").Value = 1508.4166
$ws.Range("N40").Value = -1858.4166
# row 51
$ws.Range("H51").Value = 50104010
$ws.Range("I51").Value = 203599.4
$ws.Range("J51").Value = 100004420
$ws.Range("K51").Value = 203599.4
$ws.Range("L51").Value = 100004420
$ws.Range("M51").Value = -203115.4
$ws.Range("N51").Value = -100005388
# row 112
$ws.Range("H112").Value = 3787.077
$ws.Range("J112").Value = 4089.318
$ws.Range("L112").Value = 12267.954
$ws.Range("N112").Value = -14483.954
# row 138
$ws.Range("H138").Value = 5948.622
$ws.Range("I138").Value = 100000
$ws.Range("J138").Value = 3811.0908
$ws.Range("K138").Value = 300000
$ws.Range("L138").Value = 11433.2724
$ws.Range("M138").Value = -294860
$ws.Range("N138").Value = -21713.2724

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 741.9231
$ws.Range("I2").Value = 564.8570999999999
$ws.Range("K2").Value = 564.8570999999999
$ws.Range("M2").Value = -451.8570999999999
# row 32
$ws.Range("H32").Value = 252398
$ws.Range("I32").Value = 313778.5
$ws.Range("J32").Value = 6876
$ws.Range("K32").Value = 313778.5
$ws.Range("L32").Value = 6876
$ws.Range("M32").Value = -313491.5
$ws.Range("N32").Value = -7450
# row 44
$ws.Range("H44").Value = 50023.5
# row 45
$ws.Range("H45").Value = 79480.766
$ws.Range("J45").Value = 2679.7144
$ws.Range("L45").Value = 2679.7144
$ws.Range("N45").Value = -3433.7144
# row 55
$ws.Range("H55").Value = 5000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# row 116
$ws.Range("H116").Value = 741.9231
$ws.Range("I116").Value = 564.8570999999999
$ws.Range("K116").Value = 564.8570999999999
$ws.Range("M116").Value = 1729.1429
# row 122
$ws.Range("H122").Value = 845.6667
$ws.Range("I122").Value = 845.6667
$ws.Range("K122").Value = 2537.0001
$ws.Range("M122").Value = -87.0001000000002
# row 139
$ws.Range("H139").Value = 126778
$ws.Range("J139").Value = 126778
$ws.Range("L139").Value = 126778
$ws.Range("N139").Value = -137058

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 741.9231
$ws.Range("I3").Value = 564.8570999999999
$ws.Range("K3").Value = 564.8570999999999
$ws.Range("M3").Value = -450.8570999999999
# row 20
$ws.Range("H20").Value = 940.4545000000001
$ws.Range("J20").Value = 882
$ws.Range("L20").Value = 882
$ws.Range("N20").Value = -1376
# row 105
$ws.Range("H105").Value = 56772.6
$ws.Range("J105").Value = 38333
$ws.Range("L105").Value = 38333
$ws.Range("N105").Value = -41827
# row 132
$ws.Range("H132").Value = 120000
$ws.Range("J132").Value = 120000
$ws.Range("L132").Value = 120000
$ws.Range("N132").Value = -130120

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 1009.44446
$ws.Range("I22").Value = 1001.4286
$ws.Range("K22").Value = 1001.4286
$ws.Range("M22").Value = -651.4286
# row 48
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 5000
$ws.Range("N48").Value = -5952
# row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# row 107
$ws.Range("H107").Value = 2151.1177
$ws.Range("I107").Value = 1612.1428
$ws.Range("J107").Value = 4666.3335
$ws.Range("K107").Value = 1612.1428
$ws.Range("L107").Value = 4666.3335
$ws.Range("M107").Value = 307.8571999999999
$ws.Range("N107").Value = -8506.333500000001
# row 122
$ws.Range("H122").Value = 3101.8
$ws.Range("I122").Value = 3002
$ws.Range("K122").Value = 9006
$ws.Range("M122").Value = -6556

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row 32
$ws.Range("H32").Value = 6042.5713
$ws.Range("I32").Value = 2859.8
$ws.Range("K32").Value = 8579.400000000001
$ws.Range("M32").Value = -8296.400000000001
# row 44
$ws.Range("H44").Value = 74628.57000000001
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 74628.57000000001
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 223885.71
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -224681.71
# row 55
$ws.Range("H55").Value = 4574.5386
$ws.Range("I55").Value = 3500
$ws.Range("J55").Value = 4769.909
$ws.Range("K55").Value = 10500
$ws.Range("L55").Value = 14309.727
$ws.Range("M55").Value = -10323
$ws.Range("N55").Value = -14663.727
# row 124
$ws.Range("H124").Value = 2799.4
$ws.Range("I124").Value = 3250
$ws.Range("J124").Value = 2499
$ws.Range("K124").Value = 9750
$ws.Range("L124").Value = 7497
$ws.Range("M124").Value = -4840
$ws.Range("N124").Value = -17317
# row 132
$ws.Range("H132").Value = 2147.682
$ws.Range("I132").Value = 1946.5
$ws.Range("K132").Value = 17518.5
$ws.Range("M132").Value = -14988.5
# row 139
$ws.Range("H139").Value = 3381.8333
$ws.Range("J139").Value = 5494.5
$ws.Range("L139").Value = 16483.5
$ws.Range("N139").Value = -26763.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row 63
$ws.Range("H63").Value = 61110.5
$ws.Range("I63").Value = 61110.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 61110.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -60424.5
$ws.Range("N63").ClearContents()
# row 66
$ws.Range("H66").Value = 61110.5
$ws.Range("I66").Value = 61110.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 183331.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -179899.5
$ws.Range("N66").ClearContents()
# row 97
$ws.Range("H97").Value = 1569.4
$ws.Range("I97").Value = 1214
$ws.Range("K97").Value = 1214
$ws.Range("M97").Value = -718
# row 113
$ws.Range("H113").Value = 1102.5454
$ws.Range("I113").Value = 1150
$ws.Range("K113").Value = 1150
$ws.Range("M113").Value = 1020
# row 117
$ws.Range("H117").Value = 15000
$ws.Range("J117").Value = 15000
$ws.Range("L117").Value = 15000
$ws.Range("N117").Value = -21884
# row 126
$ws.Range("H126").Value = 2946.7778
$ws.Range("I126").Value = 2766
$ws.Range("J126").Value = 3416.8
$ws.Range("K126").Value = 8298
$ws.Range("L126").Value = 10250.4
$ws.Range("M126").Value = -5828
$ws.Range("N126").Value = -15190.4
# row 132
$ws.Range("H132").Value = 987861.9399999999
$ws.Range("I132").Value = 1368.5
$ws.Range("K132").Value = 4105.5
$ws.Range("M132").Value = -1575.5
# row 134
$ws.Range("H134").Value = 42500
$ws.Range("J134").Value = 42500
$ws.Range("L134").Value = 127500
$ws.Range("N134").Value = -132570
# row 135
$ws.Range("H135").Value = 79998.336
$ws.Range("J135").Value = 79998.336
$ws.Range("L135").Value = 79998.336
$ws.Range("N135").Value = -90138.336

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 7298.8
$ws.Range("I7").Value = 3526.4
$ws.Range("J7").Value = 9813.733
$ws.Range("K7").Value = 3526.4
$ws.Range("L7").Value = 9813.733
$ws.Range("M7").Value = -3414.4
$ws.Range("N7").Value = -10037.733
# row 22
$ws.Range("H22").Value = 5766.6313
$ws.Range("I22").Value = 2621.75
$ws.Range("K22").Value = 2621.75
$ws.Range("M22").Value = -2326.75
# row 27
$ws.Range("H27").Value = 5766.6313
$ws.Range("I27").Value = 2621.75
$ws.Range("K27").Value = 2621.75
$ws.Range("M27").Value = -2514.75
# row 46
$ws.Range("H46").Value = 13248.75
$ws.Range("I46").Value = 22097.2
$ws.Range("J46").Value = 6928.4287
$ws.Range("K46").Value = 22097.2
$ws.Range("L46").Value = 6928.4287
$ws.Range("M46").Value = -21909.2
$ws.Range("N46").Value = -7304.4287
# row 55
$ws.Range("H55").Value = 1040.5405
$ws.Range("I55").Value = 1241.7059
$ws.Range("J55").Value = 869.55
$ws.Range("K55").Value = 1241.7059
$ws.Range("L55").Value = 869.55
$ws.Range("M55").Value = -1068.7059
$ws.Range("N55").Value = -1215.55
# row 101
$ws.Range("H101").Value = 35392.4
$ws.Range("J101").Value = 35392.4
$ws.Range("L101").Value = 35392.4
$ws.Range("N101").Value = -41882.4
# row 126
$ws.Range("H126").Value = 7298.8
$ws.Range("I126").Value = 3526.4
$ws.Range("J126").Value = 9813.733
$ws.Range("K126").Value = 10579.2
$ws.Range("L126").Value = 29441.199
$ws.Range("M126").Value = -8109.200000000001
$ws.Range("N126").Value = -34381.199

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row 126
$ws.Range("H126").Value = 3016.4546
$ws.Range("J126").Value = 3163.1667
$ws.Range("L126").Value = 9489.500100000001
$ws.Range("N126").Value = -14429.5001
